$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "42×85=" "17×42="
Replace-Text "47×84=" "57×93="
Replace-Text "40×51=" "49×61="
Replace-Text "50×14=" "47×28="
Replace-Text "48×65=" "63×91="
Replace-Text "33×12=" "66×93="
Replace-Text "37×62=" "12×14="
Replace-Text "80×73=" "34×52="
Replace-Text "75×22=" "50×79="
Replace-Text "72×40=" "21×32="
Replace-Text "25×94=" "92×55="
Replace-Text "46×30=" "40×43="
Replace-Text "63×36=" "30×17="
Replace-Text "61×72=" "66×79="
Replace-Text "88×84=" "95×52="
Replace-Text "96×73=" "94×85="
Replace-Text "63×68=" "69×11="
Replace-Text "47×93=" "62×76="
Replace-Text "87×32=" "67×95="
Replace-Text "53×86=" "83×50="
Replace-Text "19×62=" "48×96="
Replace-Text "66×19=" "15×16="
Replace-Text "30×92=" "26×78="
Replace-Text "15×64=" "43×78="
Replace-Text "93×84=" "84×46="
